$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J
$ws.Range("I1").Value = "Physical Activity"
$ws.Range("J1").Value = "Hands On Time"

# New column widths (closest values reachable through the COM ColumnWidth
# setter, which snaps to 1/7-character pixel increments on save; the raw
# OOXML widths these produce are as close as possible to the target
# 10.78988764044944 / 14.08988764044944)
$ws.Columns.Item(9).ColumnWidth = 10.142857142857142
$ws.Columns.Item(10).ColumnWidth = 13.428571428571429

# Fill zeros for rows 2-6 in columns I and J
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 9).Value = 0
    $ws.Cells.Item($r, 10).Value = 0
}

# Update row 7 values
$ws.Range("E7").Value = 23
$ws.Range("F7").Value = 17
$ws.Range("G7").Value = 13
$ws.Range("H7").Value = 17
$ws.Range("I7").Value = 23
$ws.Range("J7").Value = 24
